# Update the ProductLoanInput sheet:
#  - B17 (repaymentstrategy value) changes from "Mifos style" to
#    "Penalties, Fees, Interest, Principal order", taking on the same
#    cell style as the title cell B1.
#  - The active selection moves to B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

# B17 picks up the same (title-row) cell format as B1 - left/top aligned.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Activate() | Out-Null
$ws.Range("B17").Select() | Out-Null
